# Auto update stock data
# Updates the Date_1 (column A) and EBITDA (column B) values for the latest
# snapshot rows from 2026/01/13 -> 2026/01/14, refreshing EBITDA figures.
# Values must remain plain text (matching the original inlineStr cells), so
# we force a text number format before assigning and then restore the
# cell style to Normal so no stray formatting is left behind.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($range, [string]$value) {
    $range.NumberFormat = "@"
    $range.Value = $value
    $range.Style = "Normal"
}

$updates = @(
    @{ Row = 2;  Date = "2026/01/14"; Ebitda = "8.17" },
    @{ Row = 8;  Date = "2026/01/14"; Ebitda = "8.70" },
    @{ Row = 14; Date = "2026/01/14"; Ebitda = "3.18" },
    @{ Row = 20; Date = "2026/01/14"; Ebitda = "13.60" },
    @{ Row = 26; Date = "2026/01/14"; Ebitda = "11.84" },
    @{ Row = 32; Date = "2026/01/14"; Ebitda = "28.92" },
    @{ Row = 38; Date = "2026/01/14"; Ebitda = $null },
    @{ Row = 44; Date = "2026/01/14"; Ebitda = "15.07" },
    @{ Row = 50; Date = "2026/01/14"; Ebitda = "12.29" },
    @{ Row = 56; Date = "2026/01/14"; Ebitda = "34.31" },
    @{ Row = 62; Date = "2026/01/14"; Ebitda = "11.81" },
    @{ Row = 68; Date = "2026/01/14"; Ebitda = "13.35" },
    @{ Row = 74; Date = "2026/01/14"; Ebitda = "18.96" }
)

foreach ($u in $updates) {
    $r = $u.Row
    Set-TextValue $ws.Range("A$r") $u.Date
    if ($null -ne $u.Ebitda) {
        Set-TextValue $ws.Range("B$r") $u.Ebitda
    }
}

Write-Host "Updated $($updates.Count) rows"
